$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data (Nachname, Vorname, Note Exakt, Note Gerundet) per row,
# rows re-sorted + lower-cased, with new grade columns C/D populated.
$data = @(
    @("asdf",     "marlene",  5.175109999999999, 5.25),
    @("kohler",   "alina",    1,                 1),
    @("kohler",   "nina",     5.37826,            5.5),
    @("matumona", "noe",      6,                  6),
    @("matumona", "nina",     5.82464,            5.75),
    @("sarman",   "dominik",  4.78031,            4.75),
    @("zillig",   "nicolas",  3.94643,            4)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
